$d = $word.ActiveDocument

# --- Title ---
$d.Content.Find.Execute(
    "Digital Revolutions: Reshaping Societies", $false, $false, $false, $false, $false,
    $true, 1, $false, "Exploring the Realm of Art and Creativity", 2)

# --- Author name: "Sofia Perez" -> "Dr. Sarah Thompson" ---
$d.Content.Find.Execute(
    "Sofia Perez", $false, $false, $false, $false, $false,
    $true, 1, $false, "Dr. Sarah Thompson", 2)

# --- Author email user/domain ---
$d.Content.Find.Execute(
    "sperez@cambridge", $false, $false, $false, $false, $false,
    $true, 1, $false, "thomsonsarah@brookstone", 2)

# --- Body paragraph (paragraph index 4, 1-based Paragraphs(5)) full rewrite ---
$p4 = $d.Paragraphs(5).Range
$body4 = $d.Range($p4.Start, $p4.End - 1)
$body4.Text = "Art, in its myriad manifestations, transcends boundaries and captivates souls, painting a vivid tapestry of human expression. From prehistoric cave paintings to contemporary masterpieces, art reflects the collective consciousness, mirroring societal shifts and cultural paradigms. It speaks to us in a universal language, traversing time, space, and cultural divides. As we embark on a journey through the realm of art, let us unravel the enigma of its significance and explore the depths of human creativity.`v`vUnveiling the Profound Impact of Art:`v`vArt, in its profound essence, offers a window into the human experience, capturing both our triumphs and despair. It transports us to landscapes unknown, unearthing emotions we never knew we possessed. Through paintings, sculptures, music, literature, and the performing arts, we connect with diverse perspectives, fostering empathy and understanding. Art can challenge our beliefs, alter our perceptions, and ignite introspection, urging us to confront our inner demons and transcending limitations.`v`vExploring the Symphony of Creativity:`v`vCreativity, like a cascading waterfall, flows through us, an unstoppable force propelling us to create. It's a delicate dance between intuition and intellect, a harmonious interplay that transforms the mundane into the extraordinary. Creativity allows us to express ourselves and communicate our thoughts and emotions in ways words cannot capture. Irrespective of our artistic abilities, each of us possesses a creative spark, waiting to be ignited. Education can nurture this inherent creativity, guiding us to discover our unique voices and unleash our artistic potential.`v`vHarnessing Art for Societal Enrichment:`v`vArt is not merely an ornament or a pastime; it is a powerful instrument of social change. Throughout history, art has been a catalyst for awareness, provoking critical thought, and mobilizing action. Whether it be Goya's graphic depictions of the Spanish Civil War or Picasso's haunting portrayal of Guernica, art has been a clarion call to conscience, exposing injustices and advocating for transformative change. Art can amplify marginalized voices, bridge cultural divides, and foster unity. It holds the potential to inspire peace, promote understanding, and facilitate intercultural dialogue."

# --- Summary body paragraph (paragraph index 6, 1-based Paragraphs(7)) full rewrite ---
$p6 = $d.Paragraphs(7).Range
$body6 = $d.Range($p6.Start, $p6.End - 1)
$body6.Text = "Art, in its myriad manifestations, uplifts the human spirit, fosters creativity, and serves as a catalyst for social change. It challenges perceptions, expands consciousness, and connects us to our shared humanity. As we navigate an increasingly complex world, art remains an indispensable compass, guiding us towards beauty, truth, and understanding."

# --- Append a new empty paragraph at the very end of the document ---
$sel = $word.Selection
$sel.EndKey(6, 0)
$sel.TypeParagraph()

Write-Output "edit complete"
